# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" holdings sheet right before the "总计" (totals)
# sheet, populates it with the quarter's single fund holding, and updates
# the "总计" sheet with a new leading summary row for 2022-Q1 (renumbering
# the existing index column underneath it).

$wb = $excel.ActiveWorkbook

# Locate the "总计" (totals) sheet -- it is currently the last sheet -- and
# a same-shaped template sheet (any of the quarterly holdings sheets) to
# copy cell formatting from.
$totalSheet = $wb.Worksheets.Item($wb.Worksheets.Count())
$templateSheet = $wb.Worksheets.Item(3)

# A genuinely blank, never-touched cell used purely as a "no formatting"
# source for PasteSpecial(xlPasteFormats) resets later on.
$blank = $templateSheet.Range("Z100")

# Insert the new sheet immediately before the totals sheet.
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# NOTE: the worksheet collection resolves `Item(n)` positionally. Now that a
# sheet has been inserted in front of it, `$totalSheet` (captured by
# position) would resolve to the newly inserted sheet instead of "总计" if
# used as-is. Re-fetch it by name so later edits land on the right sheet.
$totalSheet = $wb.Worksheets.Item("总计")

# --- Populate the new "2022-Q1" sheet -----------------------------------

# Copy header (row 1, columns B:H) and data-row (row 2, columns A:H) cell
# formatting from the template sheet. Column A has no header cell in row 1,
# so it is deliberately excluded from that copy to avoid materialising a
# stray formatted-but-empty A1 cell.
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$templateSheet.Range("A2:H2").Copy()
$newSheet.Range("A2:H2").PasteSpecial(-4122)

# Header row.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data row. The numeric-looking text values (fund code / scale / position
# figures) must stay text, matching the source data, so their Text number
# format is applied *before* the value is assigned -- otherwise Excel
# silently reinterprets the literal as a number (dropping the leading zero
# in "001541" or the trailing zero in "2.30").
$newSheet.Range("A2").Value = 0

foreach ($addr in @("B2", "D2", "E2", "F2", "G2")) {
    $newSheet.Range($addr).NumberFormat = "@"
}

$newSheet.Range("B2").Value = "001541"
$newSheet.Range("C2").Value = "汇添富民营新动力股票"
$newSheet.Range("D2").Value = "2.30"
$newSheet.Range("E2").Value = "89.06"
$newSheet.Range("F2").Value = "3.99"
$newSheet.Range("G2").Value = "0.0918"
$newSheet.Range("H2").Value = 6

# Reset those cells' formatting back to the unstyled base look (matching
# the rest of the B2:G2 data row) now that the text format did its job --
# only the number format mattered, not a lasting style slot.
$blank.Copy()
$newSheet.Range("B2").PasteSpecial(-4122)
$newSheet.Range("D2:G2").PasteSpecial(-4122)

# --- Update the "总计" sheet ---------------------------------------------
# Insert a new row above the existing 2021-Q3 row for the 2022-Q1 summary.
$totalSheet.Rows.Item(2).Insert()

# Carry over the index-column (A) formatting from the row that used to be
# row 2 (now row 3, still bearing the bold/bordered index style) onto the
# freshly inserted row.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.09

# Row-insert copied the header row's bold formatting down onto B2:D2;
# reset those back to the unstyled look the other data rows use.
$blank.Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122)

# Renumber the remaining index (A) column values: 1, 2, 3
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
